$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New data rows to append to the table
$ws.Range("A11").Value = 109
$ws.Range("B11").Value = "SESSION_FAILED"
$ws.Range("C11").Value = "Nebolo možné utvoriť reláciu medzi databázou a klientom"
$ws.Range("D11").Value = "Skontrolovať prihlasovacie údaje"

$ws.Range("A12").Value = 110
$ws.Range("B12").Value = "REFFERENCES_RETRIEVAL_FAILURE"
$ws.Range("D12").Value = "Kontaktovať administrátora"
$ws.Range("C12").Value = "Pri sťahovaní dát z jednej z tabuliek: Path,Actor,Board,HDV,Software nastala chyba"

# Resize the table to include the new rows
$tbl = $ws.ListObjects.Item("Table1")
$tbl.Resize($ws.Range("A1:D12"))

# Adjust column widths to fit new (longer) content (bestFit recalculated by Excel)
$ws.Columns.Item(2).ColumnWidth = 31.2
$ws.Columns.Item(3).ColumnWidth = 74.2

# Update the selection to match the post-edit state
$ws.Range("C14").Select()
